$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Golang Developer", "https://www.dice.com/job-detail/88326c4f-7e6e-420e-a59d-999b7caff3d7", "Fremont, California", "Contract, Third Party", "Depends on Experience", "TechVirtue LLC"),
    @("Lead Golang Developer", "https://www.dice.com/job-detail/83430e9f-d44c-4265-a77d-726d4701dfcc", "Richmond, Virginia", "Contract, Third Party", "Depends on Experience", "NasTech Global, Inc.")
)

$startRow = 108
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$c]
    }
}
